$d = $word.ActiveDocument

# The document contains three "<id>...</id>" markers (p124v_1, p124v_2,
# p124v_3) that were each split across three separate runs:
#   run1 "<id>"  (Courier New)  + run2 "p124v_N" (Arial) + run3 "</id>" (Courier New)
# Newly downloaded tc/tcn/tl content supplies these ids as a single run,
# so collapse each split marker back into one run: "<id>p124v_N</id>".
# Using Find/Replace (rather than touching XML runs directly) lets Word
# merge the matched runs and carry over the formatting of the first
# character of the match, exactly as Word does interactively.

$ids = @("p124v_1", "p124v_2", "p124v_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $old, 2) | Out-Null
}
